$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used below to strip any incidental style (e.g. the
# "quote prefix" / "stored as text" style Excel applies when a value looks
# like a number) back down to the plain default style, matching a cell that
# was never specially formatted.
$xlPasteFormats = -4122

# --- 1. Insert a new row at position 17 -----------------------------------
# Shifts old row 17 ("Commonly Used React Hooks...") down to row 18, and old
# row 18 ("...ssh connection settings...") down to row 19.
$ws.Rows.Item(17).Insert()

# --- 2. Populate the newly inserted row 17 with the new post's data -------
$ws.Cells.Item(17, 1).Value = "Draft"
$ws.Cells.Item(17, 2).Value = "Article"

# C17 ("id") is textual ("17") in the source data, like the other id cells
# in this column. Writing the bare digit string makes Excel auto-coerce it
# to a number, so force text storage, then strip the resulting "quote
# prefix" style back to the default (unstyled) look.
$c17 = $ws.Cells.Item(17, 3)
$c17.Value = "'17"
$ws.Cells.Item(17, 1).Copy()
$c17.PasteSpecial($xlPasteFormats)

$ws.Cells.Item(17, 4).Value = "The Human Factors of Formal Methods"
$ws.Cells.Item(17, 5).Value = "As formal methods improve in expressiveness and power, they create new opportunities for non-expert adoption. In principle, formal tools are now powerful enough to enable developers to scalably validate realistic systems artifacts without extensive formal training."
$ws.Cells.Item(17, 6).Value = "Human_Factors_Formal_Methods"

# G17 ("date") - reuse the same date-formatted style as the other date
# cells in this column (e.g. the one now at G18) so the serial number
# renders the same way (yyyy-mm-dd).
$g17 = $ws.Cells.Item(17, 7)
$g17.Value = 45399
$ws.Cells.Item(18, 7).Copy()
$g17.PasteSpecial($xlPasteFormats)

$ws.Cells.Item(17, 8).Value = "Shoto Morisaki"
$ws.Cells.Item(17, 9).Value = "sample/sample7.jpg"
$ws.Cells.Item(17, 10).Value = "Article"
$ws.Cells.Item(17, 11).Value = "['Talk']"

# L17:P17 (relatedDoc1-5) are blank in the source data, but still exist as
# empty text cells (same as every other row in this sheet). A plain
# Value = "" does not persist an empty cell, so write a lone apostrophe
# (an empty quote-prefixed string) to force the cell to exist as text, then
# strip the resulting style back to the plain default.
foreach ($col in 12..16) {
    $cell = $ws.Cells.Item(17, $col)
    $cell.Value = "'"
    $ws.Cells.Item(17, 1).Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

$ws.Cells.Item(17, 17).Value = "https://raw.githubusercontent.com/shoot649854/IMG_DB/main/profile.webp"
$ws.Cells.Item(17, 18).Value = "Computer Science - University of California Santa Cruz | Intern - LiNK"
$ws.Cells.Item(17, 19).Value = "Shoto Morisaki"
$ws.Cells.Item(17, 20).Value = "[{'icon': 'fab fa-linkedin', 'url': 'https://www.linkedin.com/in/shoto-morisaki-93b0a71bb/'}, {'icon': 'fab fa-github', 'url': 'https://github.com/shoot649854/'}, {'icon': 'fab fa-portfolio', 'url': 'https://portfolio-shoto.vercel.app/'}]"

$ws.Application.CutCopyMode = $false

# --- 3. Update the tags in row 4 (K4) --------------------------------------
$ws.Cells.Item(4, 11).Value = "['AWS', 'Docker', 'Python', 'Flask', 'MUI', 'TypeScript', 'React']"

# --- 4. Replace the old GitHub "blob" profile image URL with the new ------
#        raw.githubusercontent URL across every data row, column Q (17)
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 17)
    if ($cell.Value2 -eq "https://github.com/shoot649854/IMG_DB/blob/main/profile.webp") {
        $cell.Value = "https://raw.githubusercontent.com/shoot649854/IMG_DB/main/profile.webp"
    }
}
